# Test Data Comment column updated for login and reset password
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Login_Page sheet
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Login_Page")

# New "Comment" header + make the whole header row bold (A1:B1 already bold)
$ws1.Range("C1").Value = "Comment"
$ws1.Range("A1:C1").Font.Bold = $true

# New column C comments for existing + new rows
$ws1.Range("C2").Value = "Valid Credential"
$ws1.Range("C3").Value = "Invalid Credentials"
$ws1.Range("C4").Value = "Blank User Name,Valid Passwprd"
$ws1.Range("C5").Value = "Valid User name , Blank Password"

# New row 6: extra test case
$ws1.Range("A6").Value = "LmsAutomation"
$ws1.Range("B6").Value = "sdffsdsd"
$ws1.Range("C6").Value = "Valid Username,Invalid pw"

# New row 7: extra test case
$ws1.Range("A7").Value = "dfdssss"
$ws1.Range("B7").Value = "abc@123"
$ws1.Range("C7").Value = "InValid Username,Valid pw"

# New row 8: comment only
$ws1.Range("C8").Value = "Empty User name and password"

# Widen the new comment column and move the selection like the source file
$ws1.Columns(3).ColumnWidth = 30.72
$ws1.Range("C18").Select() | Out-Null

# ---------------------------------------------------------------------------
# Reset_Password sheet
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Reset_Password")

# New "Comment" header + bold header row (A1:B1 previously unformatted)
$ws2.Range("C1").Value = "Comment"
$ws2.Range("A1:C1").Font.Bold = $true

# New column C comments
$ws2.Range("C2").Value = "Valid details"
$ws2.Range("C2").Font.Name = "Consolas"
$ws2.Range("C2").Font.Size = 12
$ws2.Range("C2").Font.Color = 0
$ws2.Range("C2").VerticalAlignment = -4108
$ws2.Rows(2).RowHeight = 15.6

$ws2.Range("C3").Value = "Invalid details"
$ws2.Range("C4").Value = "Empty details"
$ws2.Range("C5").Value = "Mismatch values"

# Widen the new comment column and move the selection like the source file
$ws2.Columns(3).ColumnWidth = 24.17
$ws2.Range("C16").Select() | Out-Null
